# edit.ps1
# Adds a new "2022-Q3" quarter sheet (with its fund-holding detail data) to the
# workbook, right after the "总计" (Total) summary sheet, and inserts the
# corresponding summary row at the top of the "总计" sheet's data table.
#
# All of the other existing quarter sheets (2022-Q2, 2022-Q1, 2021-Q4,
# 2021-Q3, 2021-Q1, 2020-Q4) keep their names/content unchanged - they are
# simply pushed one position later in the sheet order by the insertion.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" (Total) sheet: insert a new row under the header for 2022-Q3,
#    pushing the existing quarter rows down by one.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

$total.Rows.Item(2).Insert()

# `Insert()` copies the format of the row above (the bold header row) onto
# the whole new row; the data columns (B:D) should stay plain like every
# other data row, so strip that back off before writing values.
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 12
$total.Range("D2").Value = 5.47

# Match the look of the other index cells in column A (bold / bordered /
# centered) for the newly inserted row's index cell.
$total.Range("A2").Font.Bold = $true
$total.Range("A2").HorizontalAlignment = -4108
$total.Range("A2").VerticalAlignment = -4160
$total.Range("A2").Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 2. Add the new "2022-Q3" sheet right after "总计".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}

$rows = @(
    @(0,  "011363", "南方兴润价值一年持有期混合A", "67.01", "65.71", "3.87", "2.5933", 4),
    @(1,  "202003", "南方绩优成长混合A",           "47.57", "73.75", "4.06", "1.9313", 1),
    @(2,  "011364", "南方兴润价值一年持有期混合C", "18.36", "65.71", "3.87", "0.7105", 4),
    @(3,  "014202", "天弘中证1000指数增强C",       "3.69",  "94.06", "1.60", "0.0590", 3),
    @(4,  "014201", "天弘中证1000指数增强A",       "3.68",  "94.06", "1.60", "0.0589", 3),
    @(5,  "015784", "中信建投中证1000指数增强A",   "8.10",  "92.20", "0.66", "0.0535", 6),
    @(6,  "015785", "中信建投中证1000指数增强C",   "3.32",  "92.20", "0.66", "0.0219", 6),
    @(7,  "002872", "华夏智胜价值成长股票C",       "2.13",  "93.39", "0.93", "0.0198", 4),
    @(8,  "006540", "南方绩优成长混合C",           "0.27",  "73.75", "4.06", "0.0110", 1),
    @(9,  "002871", "华夏智胜价值成长股票A",       "0.86",  "93.39", "0.93", "0.0080", 4),
    @(10, "005429", "渤海汇金睿选混合A",           "0.13",  "30.88", "2.15", "0.0028", 2),
    @(11, "005430", "渤海汇金睿选混合C",           "0.01",  "30.88", "2.15", "0.0002", 2)
)

# Columns B (fund code), D (fund scale), E (stock position), F (position
# ratio), G (market value held) all carry leading-zero codes / fixed-decimal
# figures stored as plain TEXT in the source data (not numbers), so force a
# text number-format before writing them - otherwise e.g. "011363" would be
# coerced to the number 11363.
$q3.Range("B2:B13").NumberFormat = "@"
$q3.Range("D2:G13").NumberFormat = "@"

$r = 2
foreach ($row in $rows) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).Value = $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r++
}

# Header row + index column formatting, matching the other quarter sheets
# (bold / bordered / centered).
$q3Header = $q3.Range("B1:H1")
$q3Header.Font.Bold = $true
$q3Header.HorizontalAlignment = -4108
$q3Header.VerticalAlignment = -4160
$q3Header.Borders.LineStyle = 1

$q3Index = $q3.Range("A2:A13")
$q3Index.Font.Bold = $true
$q3Index.HorizontalAlignment = -4108
$q3Index.VerticalAlignment = -4160
$q3Index.Borders.LineStyle = 1
